$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the numeric values in columns B:E (rows 2-13) to the nearest integer,
# matching the behaviour of writing the data to disk as integers.
$range = $ws.Range("B2:E13")
foreach ($cell in $range.Cells) {
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value2 = $excel.WorksheetFunction.Round([double]$val, 0)
    }
}
